## Applies the "List of Inputs" sheet updates described in the commit diff:
##  - appends a follow-up note to the "how to make a variable?" comment (C5)
##  - fills in notes / hyperlink / data-years for the TSO/DSO financial
##    statements and energy-demand-forecasting rows (C9, C10, D10, C11, C12)
##  - narrows column B and widens column C
##  - moves the saved cell selection from C8 to C6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List of Inputs")

# --- Row 5: extend the existing comment in column C -----------------------
$ws.Range("C5").Value = "how to make a variable? Maybe use NLP? \n Whic h links to use? How to interpret the data?"

# --- Row 9: new comment about TSO 10-year infrastructure plans ------------
$ws.Range("C9").Value = "found one publication"

# --- Row 10: TSO/DSO financial statements - comment, hyperlink, years -----
# Match column C's existing wrap/centre formatting before the Hyperlinks.Add
# call switches the cell to the built-in hyperlink style.
$ws.Range("C10").HorizontalAlignment = $ws.Range("C9").HorizontalAlignment
$ws.Range("C10").VerticalAlignment = $ws.Range("C9").VerticalAlignment
$ws.Range("C10").WrapText = $ws.Range("C9").WrapText

$ws.Hyperlinks.Add($ws.Range("C10"), "https://sse.com/investors/reportsandresults/reports/", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "https://sse.com/investors/reportsandresults/reports/")
$ws.Range("C10").Value = "There are financial statement, but what exactly to look for. /n https://sse.com/investors/reportsandresults/reports/"

# Restore the plain (non-hyperlink) look used by the rest of column C.
$ws.Range("C10").HorizontalAlignment = $ws.Range("C9").HorizontalAlignment
$ws.Range("C10").VerticalAlignment = $ws.Range("C9").VerticalAlignment
$ws.Range("C10").WrapText = $ws.Range("C9").WrapText
$ws.Range("C10").Font.Underline = $false
$ws.Range("C10").Font.Color = $ws.Range("C9").Font.Color
$ws.Range("C10").Font.Name = $ws.Range("C9").Font.Name
$ws.Range("C10").Font.Size = $ws.Range("C9").Font.Size

$ws.Range("D10").Value = "19 years"

# --- Row 11: National energy demand forecasting ----------------------------
$ws.Range("C11").Value = "I could not find the repos for that"

# --- Row 12: National energy and Climate plans -----------------------------
$ws.Range("C12").Value = "?"

# --- Column widths: narrower Input Name column, wider comments column -----
$ws.Columns.Item(2).ColumnWidth = 47.8
$ws.Columns.Item(3).ColumnWidth = 32.8

# --- Saved selection moves from C8 to C6 -----------------------------------
$ws.Range("C6").Select()

Write-Host "Applied List of Inputs updates"
